# Auto-generated edit script applying the Marilith_Profits.xlsx diff
# Updates leve-crafting profit calculations across ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1929
$ws.Range("I43").Value = 1893.5
$ws.Range("K43").Value = 1893.5
$ws.Range("M43").Value = -1824.5
$ws.Range("H98").Value = 3241.2778
$ws.Range("I98").Value = 3064
$ws.Range("J98").Value = 3595.8333
$ws.Range("K98").Value = 3064
$ws.Range("L98").Value = 3595.8333
$ws.Range("M98").Value = -1566
$ws.Range("N98").Value = -6591.8333
$ws.Range("H100").Value = 3787
$ws.Range("I100").Value = 4047.3333
$ws.Range("J100").Value = 3006
$ws.Range("K100").Value = 4047.3333
$ws.Range("L100").Value = 3006
$ws.Range("M100").Value = -3506.3333
$ws.Range("N100").Value = -4088
$ws.Range("H111").Value = 1699.75
$ws.Range("I111").Value = 1266.3334
$ws.Range("K111").Value = 3799.0002
$ws.Range("M111").Value = -732.0001999999999
$ws.Range("H112").Value = 2237.375
$ws.Range("J112").Value = 2342.7856
$ws.Range("L112").Value = 7028.3568
$ws.Range("N112").Value = -9244.356800000001
$ws.Range("H122").Value = 3241.2778
$ws.Range("I122").Value = 3064
$ws.Range("J122").Value = 3595.8333
$ws.Range("K122").Value = 9192
$ws.Range("L122").Value = 10787.4999
$ws.Range("M122").Value = -6742
$ws.Range("N122").Value = -15687.4999
$ws.Range("H137").Value = 1352
$ws.Range("I137").Value = 1211.3158
$ws.Range("J137").Value = 2020.25
$ws.Range("K137").Value = 3633.9474
$ws.Range("L137").Value = 6060.75
$ws.Range("M137").Value = -1083.9474
$ws.Range("N137").Value = -11160.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6101.0884
$ws.Range("I32").Value = 4770.8184
$ws.Range("K32").Value = 4770.8184
$ws.Range("M32").Value = -4483.8184
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35976
$ws.Range("H61").Value = 1498.6666
$ws.Range("I61").Value = 1498.6666
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1498.6666
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1286.6666
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2231
$ws.Range("I74").Value = 2231
$ws.Range("K74").Value = 2231
$ws.Range("M74").Value = -1357
$ws.Range("H77").Value = 2231
$ws.Range("I77").Value = 2231
$ws.Range("K77").Value = 11155
$ws.Range("M77").Value = -6787
$ws.Range("H122").Value = 2419.1177
$ws.Range("I122").Value = 2080.4285
$ws.Range("K122").Value = 6241.2855
$ws.Range("M122").Value = -3791.2855
$ws.Range("H132").Value = 1748.579
$ws.Range("I132").Value = 1830.6
$ws.Range("J132").Value = 1441
$ws.Range("K132").Value = 5491.799999999999
$ws.Range("L132").Value = 4323
$ws.Range("M132").Value = -2961.799999999999
$ws.Range("N132").Value = -9383
$ws.Range("H136").Value = 1498.6666
$ws.Range("I136").Value = 1498.6666
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4495.9998
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1945.9998
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1543.5714
$ws.Range("I20").Value = 761.2222
$ws.Range("K20").Value = 761.2222
$ws.Range("M20").Value = -514.2222
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H75").Value = 27671.5
$ws.Range("J75").Value = 130000
$ws.Range("L75").Value = 130000
$ws.Range("N75").Value = -131872
$ws.Range("H78").Value = 27671.5
$ws.Range("J78").Value = 130000
$ws.Range("L78").Value = 390000
$ws.Range("N78").Value = -399360
$ws.Range("H86").Value = 1657.8125
$ws.Range("I86").Value = 1713.2
$ws.Range("J86").Value = 1565.5
$ws.Range("K86").Value = 1713.2
$ws.Range("L86").Value = 1565.5
$ws.Range("M86").Value = -590.2
$ws.Range("N86").Value = -3811.5
$ws.Range("H89").Value = 1657.8125
$ws.Range("I89").Value = 1713.2
$ws.Range("J89").Value = 1565.5
$ws.Range("K89").Value = 8566
$ws.Range("L89").Value = 7827.5
$ws.Range("M89").Value = -2950
$ws.Range("N89").Value = -19059.5
$ws.Range("H134").Value = 7910.1
$ws.Range("I134").Value = 7910.1
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 23730.3
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -21195.3
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3454.4666
$ws.Range("I16").Value = 1899.1818
$ws.Range("K16").Value = 1899.1818
$ws.Range("M16").Value = -1612.1818
$ws.Range("H31").Value = 2342.611
$ws.Range("I31").Value = 2409.125
$ws.Range("J31").Value = 2289.4
$ws.Range("K31").Value = 2409.125
$ws.Range("L31").Value = 2289.4
$ws.Range("M31").Value = -2114.125
$ws.Range("N31").Value = -2879.4
$ws.Range("H34").Value = 2342.611
$ws.Range("I34").Value = 2409.125
$ws.Range("J34").Value = 2289.4
$ws.Range("K34").Value = 2409.125
$ws.Range("L34").Value = 2289.4
$ws.Range("M34").Value = -2207.125
$ws.Range("N34").Value = -2693.4
$ws.Range("H113").Value = 3454.4666
$ws.Range("I113").Value = 1899.1818
$ws.Range("K113").Value = 1899.1818
$ws.Range("M113").Value = 270.8181999999999
$ws.Range("H134").Value = 2937.6428
$ws.Range("I134").Value = 1456.2
$ws.Range("K134").Value = 4368.6
$ws.Range("M134").Value = -1833.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7818227.5
$ws.Range("I122").Value = 13895760
$ws.Range("J122").Value = 4257
$ws.Range("K122").Value = 41687280
$ws.Range("L122").Value = 12771
$ws.Range("M122").Value = -41684830
$ws.Range("N122").Value = -17671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3585.4285
$ws.Range("I46").Value = 550
$ws.Range("K46").Value = 550
$ws.Range("M46").Value = -362
$ws.Range("H61").Value = 7216.857
$ws.Range("I61").Value = 6869.8335
$ws.Range("K61").Value = 6869.8335
$ws.Range("M61").Value = -6667.8335
$ws.Range("H68").Value = 5812.5
$ws.Range("J68").Value = 6125
$ws.Range("L68").Value = 6125
$ws.Range("N68").Value = -7623
$ws.Range("H71").Value = 5812.5
$ws.Range("J71").Value = 6125
$ws.Range("L71").Value = 30625
$ws.Range("N71").Value = -38113
$ws.Range("H113").Value = 7216.857
$ws.Range("I113").Value = 6869.8335
$ws.Range("K113").Value = 6869.8335
$ws.Range("M113").Value = -4699.8335
$ws.Range("H122").Value = 4063.75
$ws.Range("I122").Value = 3501.4285
$ws.Range("K122").Value = 10504.2855
$ws.Range("M122").Value = -8054.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 60000
$ws.Range("J70").Value = 60000
$ws.Range("L70").Value = 60000
$ws.Range("N70").Value = -60630
$ws.Range("H73").Value = 60000
$ws.Range("J73").Value = 60000
$ws.Range("L73").Value = 60000
$ws.Range("N73").Value = -62184
$ws.Range("H107").Value = 1436.6
$ws.Range("I107").Value = 1436.6
$ws.Range("K107").Value = 4309.799999999999
$ws.Range("M107").Value = -2389.799999999999
$ws.Range("H126").Value = 3424.4443
$ws.Range("J126").Value = 3617.5
$ws.Range("L126").Value = 10852.5
$ws.Range("N126").Value = -15792.5
$ws.Range("H132").Value = 1053.25
$ws.Range("I132").Value = 1104.4445
$ws.Range("J132").Value = 899.6667
$ws.Range("K132").Value = 3313.3335
$ws.Range("L132").Value = 2699.0001
$ws.Range("M132").Value = -783.3335000000002
$ws.Range("N132").Value = -7759.0001
